$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(232, 8, '30mm SL', 2, 'A1', 'TCMM164'),
    @(233, 8, '26mm SL', 2, 'A2', 'TCMM165'),
    @(234, 8, '28mm SL', 2, 'A3', 'TCMM166'),
    @(235, 8, '30mm SL', 2, 'A4', 'TCMM167'),
    @(236, 8, '29mm SL', 2, 'A5', 'TCMM168'),
    @(237, 8, '29mm SL', 2, 'A6', 'TCMM169'),
    @(238, 8, '31mm SL', 2, 'A7', 'TCMM170'),
    @(239, 8, '29mm SL', 2, 'A8', 'TCMM171'),
    @(240, 8, '29mm SL', 2, 'A9', 'TCMM172'),
    @(241, 8, '27mm SL', 2, 'A10', 'TCMM173'),
    @(242, 8, '28mm SL', 2, 'A11', 'TCMM174'),
    @(243, 8, '29mm SL', 2, 'A12', 'TCMM175'),
    @(244, 8, '32mm SL', 2, 'B1', 'TCMM176'),
    @(245, 8, '28mm SL', 2, 'B2', 'TCMM177'),
    @(246, 8, '29mm SL', 2, 'B3', 'TCMM178'),
    @(247, 8, '33mm SL', 1, 'B4', 'TCMM179'),
    @(248, 8, '29mm SL', 2, 'B5', 'TCMM180'),
    @(249, 8, '26mm SL', 2, 'B6', 'TCMM181'),
    @(250, 8, '27mm SL', 1, 'B7', 'TCMM182'),
    @(251, 8, '26mm SL', 1, 'B8', 'TCMM183'),
    @(252, 8, '33mm SL', 2, 'B9', 'TCMM184'),
    @(253, 8, '28mm SL', 2, 'B10', 'TCMM185'),
    @(254, 8, '29mm SL', 2, 'B11', 'TCMM186'),
    @(255, 8, '29mm SL', 2, 'B12', 'TCMM187'),
    @(256, 8, '29mm SL', 2, 'C1', 'TCMM188'),
    @(257, 8, '30mm SL', 2, 'C2', 'TCMM189'),
    @(258, 8, '26mm SL', 2, 'C3', 'TCMM190'),
    @(259, 8, '28mm SL', 1, 'C4', 'TCMM191'),
    @(260, 8, '28mm SL', 2, 'C5', 'TCMM192'),
    @(261, 8, '30mm SL', 2, 'C6', 'TCMM193'),
    @(262, 8, '28mm SL', 2, 'C7', 'TCMM194'),
    @(263, 8, '26mm SL', 2, 'C8', 'TCMM195'),
    @(264, 8, '30mm SL', 2, 'C9', 'TCMM196'),
    @(265, 8, '28mm SL', 2, 'C10', 'TCMM197'),
    @(266, 8, '32mm SL', 2, 'C11', 'TCMM198'),
    @(267, 8, '30mm SL', 2, 'C12', 'TCMM199'),
    @(268, 8, '32mm SL', 2, 'D1', 'TCMM200'),
    @(269, 8, '29mm SL', 2, 'D2', 'TCMM201'),
    @(270, 8, '33mm SL', 2, 'D3', 'TCMM202'),
    @(271, 8, '31mm SL', 2, 'D4', 'TCMM203'),
    @(272, 8, '26mm SL', 2, 'D5', 'TCMM204'),
    @(273, 8, '28mm SL', 2, 'D6', 'TCMM205'),
    @(274, 8, '32mm SL', 2, 'D7', 'TCMM206'),
    @(275, 8, '28mm SL', 2, 'D8', 'TCMM207'),
    @(276, 8, '29mm SL', 1, 'D9', 'TCMM208'),
    @(277, 8, '29mm SL', 2, 'D10', 'TCMM209'),
    @(278, 8, '27mm SL', 2, 'D11', 'TCMM210'),
    @(279, 8, '25mm SL', 2, 'D12', 'TCMM211'),
    @(280, 8, '30mm SL', 2, 'E1', 'TCMM212'),
    @(281, 8, '33mm SL', 2, 'E2', 'TCMM213'),
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = 'Maurolicus muelleri'
    $ws.Cells.Item($row, 2).Font.Italic = $true
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = 'TC5'
    $ws.Cells.Item($row, 7).Value = $r[5]
}

$excel.ActiveWindow.Zoom = 53
$ws.Range("D282").Select()
"done"